# Apply updated crypto price/volume figures (and re-ranked rows 10-15)
# to match the refreshed data snapshot, per commit "Updated symbol list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'255.93"
$ws.Range("E2").Value = "'4.28%"

# Row 3
$ws.Range("D3").Value = "'28.00"
$ws.Range("E3").Value = "'-4.65%"

# Row 4
$ws.Range("D4").Value = "'5.353"
$ws.Range("E4").Value = "'3.85%"

# Row 5
$ws.Range("D5").Value = "'0.05827"
$ws.Range("E5").Value = "'0.86%"

# Row 6
$ws.Range("D6").Value = "'6.710"
$ws.Range("E6").Value = "'1.24%"

# Row 7
$ws.Range("E7").Value = "'2.49%"

# Row 8
$ws.Range("D8").Value = "'0.8722"
$ws.Range("E8").Value = "'1.82%"

# Row 9
$ws.Range("D9").Value = "'0.9166"
$ws.Range("E9").Value = "'7.10%"

# Row 10
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01052"
$ws.Range("E10").Value = "'1,647.78%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1419"
$ws.Range("E11").Value = "'4.01%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07189"
$ws.Range("E12").Value = "'2.28%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03180"
$ws.Range("E13").Value = "'4.28%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09245"
$ws.Range("E14").Value = "'-1.32%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001553"
$ws.Range("E15").Value = "'1.44%"

# Row 16
$ws.Range("D16").Value = "'0.005958"
$ws.Range("E16").Value = "'-1.13%"

# Row 17
$ws.Range("E17").Value = "'0.46%"

# Row 18
$ws.Range("D18").Value = "'2.274"
$ws.Range("E18").Value = "'5.24%"

# Row 19
$ws.Range("D19").Value = "'0.3166"
$ws.Range("E19").Value = "'-1.16%"

# Row 20
$ws.Range("D20").Value = "'0.03449"
$ws.Range("E20").Value = "'4.26%"

# Row 21
$ws.Range("E21").Value = "'2.22%"

# Row 22
$ws.Range("D22").Value = "'3.525"
$ws.Range("E22").Value = "'6.25%"

# Row 23
$ws.Range("D23").Value = "'0.04156"
$ws.Range("E23").Value = "'0.68%"

# Row 24
$ws.Range("E24").Value = "'-1.56%"

# Row 25
$ws.Range("D25").Value = "'0.001227"
$ws.Range("E25").Value = "'0.10%"

# Row 26
$ws.Range("D26").Value = "'0.004876"
$ws.Range("E26").Value = "'18.02%"

# Row 27
$ws.Range("E27").Value = "'-0.89%"

# Row 28
$ws.Range("E28").Value = "'0.63%"

# Row 40
$ws.Range("D40").Value = "'0.03851"

# Row 41
$ws.Range("E41").Value = "'-2.56%"

# Row 42
$ws.Range("D42").Value = "'0.1100"
$ws.Range("E42").Value = "'3.00%"

# Row 43
$ws.Range("D43").Value = "'0.002198"
$ws.Range("E43").Value = "'-0.06%"

# Row 44
$ws.Range("D44").Value = "'0.009913"
$ws.Range("E44").Value = "'8.32%"

# Row 45
$ws.Range("D45").Value = "'0.00005290"
$ws.Range("E45").Value = "'-0.07%"

# Row 47
$ws.Range("D47").Value = "'0.09992"
$ws.Range("E47").Value = "'72.30%"

# Row 48
$ws.Range("D48").Value = "'0.002128"
$ws.Range("E48").Value = "'-2.02%"

Write-Host "Applied crypto snapshot update."
